$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.152.05'
$ws.Range('E2').Value = '  +3.73%  '
$ws.Range('D3').Value = '1.894.63'
$ws.Range('E3').Value = '  +3.82%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9978'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9979'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4982'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.80'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2945'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06655'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.20%  '
$ws.Range('D11').Value = '1.892.66'
$ws.Range('E11').Value = '  +3.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.99'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07197'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6777'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '85.82'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.851'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.30%  '
$ws.Range('D17').Value = '30.123.57'
$ws.Range('E17').Value = '  +3.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008006'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9966'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('E20').Value = '  +6.10%  '
$ws.Range('D21').Value = '2.136.97'
$ws.Range('E21').Value = '  +3.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9975'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.776'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.654'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.163'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '133.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.946'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.372'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.220'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08733'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.950'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05135'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.118'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7075'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.664'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.773'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.229'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9422'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('E41').Value = '  +4.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.069'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9966'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4220'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.18'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.490'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.37%  '
$ws.Range('E47').Value = '  +3.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05725'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('E49').Value = '  +3.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.283'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3741'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.23%  '
